$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 data
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "Binary S 13"
$ws.Range("D5").Value = "Find a peak element"
$ws.Range("E5").Value = "https://www.scaler.com/academy/mentee-dashboard/class/30364/assignment/problems/4132/?navref=cl_pb_nv_tb"

# Add hyperlink on E5
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.scaler.com/academy/mentee-dashboard/class/30364/assignment/problems/4132/?navref=cl_pb_nv_tb", "", "", "https://www.scaler.com/academy/mentee-dashboard/class/30364/assignment/problems/4132/?navref=cl_pb_nv_tb")

# Apply styles consistent with rows above (row height, cell style, alignment, wrap text)
$ws.Range("B5:D5").Style = "Normal"
$ws.Range("B5:D5").HorizontalAlignment = -4131  # xlLeft
$ws.Range("B5:D5").VerticalAlignment = -4160    # xlTop

$ws.Rows.Item(5).RowHeight = 72

# Update selection
$ws.Range("F5").Select()
